$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve original (default) style of the D2:D51 price column, then force
# it to Text format so that numeric-looking strings (e.g. "1.000", "5.399")
# round-trip as literal text instead of being auto-coerced to numbers.
$origStyle = $ws.Range("D2").Style
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.280.61'
$ws.Range("E2").Value = '  -0.75%  '

$ws.Range("D3").Value = '1.887.66'
$ws.Range("E3").Value = '  -1.27%  '

$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Value = '238.60'
$ws.Range("E5").Value = '  -0.42%  '

$ws.Range("D6").Value = '1.000'

$ws.Range("D7").Value = '0.4686'
$ws.Range("E7").Value = '  -1.68%  '

$ws.Range("D8").Value = '0.2861'
$ws.Range("E8").Value = '  +0.52%  '

$ws.Range("D9").Value = '0.06610'
$ws.Range("E9").Value = '  -1.03%  '

$ws.Range("D10").Value = '20.09'
$ws.Range("E10").Value = '  +7.31%  '

$ws.Range("D11").Value = '0.07788'
$ws.Range("E11").Value = '  +1.34%  '

$ws.Range("D12").Value = '98.37'
$ws.Range("E12").Value = '  -2.98%  '

$ws.Range("D13").Value = '1.898.38'
$ws.Range("E13").Value = '  -0.76%  '

$ws.Range("D14").Value = '5.134'
$ws.Range("E14").Value = '  -1.87%  '

$ws.Range("D15").Value = '0.6786'
$ws.Range("E15").Value = '  +1.25%  '

$ws.Range("D16").Value = '286.16'
$ws.Range("E16").Value = '  +12.21%  '

$ws.Range("D17").Value = '30.282.08'
$ws.Range("E17").Value = '  -0.78%  '

$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.04%  '

$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '2.141.83'
$ws.Range("E19").Value = '  -1.03%  '

$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").Value = '12.65'
$ws.Range("E20").Value = '  -0.07%  '

$ws.Range("D21").Value = '5.399'
$ws.Range("E21").Value = '  +0.28%  '

$ws.Range("D22").Value = '0.000007314'
$ws.Range("E22").Value = '  -2.15%  '

$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  -0.08%  '

$ws.Range("D24").Value = '6.208'
$ws.Range("E24").Value = '  -1.42%  '

$ws.Range("D25").Value = '9.471'
$ws.Range("E25").Value = '  +1.44%  '

$ws.Range("D26").Value = '165.98'
$ws.Range("E26").Value = '  -1.46%  '

$ws.Range("D27").Value = '19.32'
$ws.Range("E27").Value = '  +1.61%  '

$ws.Range("D28").Value = '2.004'
$ws.Range("E28").Value = '  -3.12%  '

$ws.Range("D29").Value = '1.378'
$ws.Range("E29").Value = '  -0.09%  '

$ws.Range("D30").Value = '0.09740'
$ws.Range("E30").Value = '  -3.46%  '

$ws.Range("D31").Value = '4.475'
$ws.Range("E31").Value = '  -5.17%  '

$ws.Range("D32").Value = '1.486'
$ws.Range("E32").Value = '  -1.79%  '

$ws.Range("D33").Value = '4.173'
$ws.Range("E33").Value = '  -1.71%  '

$ws.Range("D34").Value = '0.04721'
$ws.Range("E34").Value = '  +0.06%  '

$ws.Range("D35").Value = '0.7136'
$ws.Range("E35").Value = '  -1.98%  '

$ws.Range("E36").Value = '  -0.71%  '

$ws.Range("D37").Value = '2.709'
$ws.Range("E37").Value = '  +0.15%  '

$ws.Range("D38").Value = '0.01886'
$ws.Range("E38").Value = '  -1.20%  '

$ws.Range("D39").Value = '6.666'
$ws.Range("E39").Value = '  +6.95%  '

$ws.Range("D40").Value = '2.531'
$ws.Range("E40").Value = '  -3.23%  '

$ws.Range("D41").Value = '72.83'
$ws.Range("E41").Value = '  -2.85%  '

$ws.Range("D42").Value = '1.979'
$ws.Range("E42").Value = '  +0.60%  '

$ws.Range("D43").Value = '0.8729'
$ws.Range("E43").Value = '  +1.32%  '

$ws.Range("D44").Value = '104.28'
$ws.Range("E44").Value = '  -0.71%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").Value = '0.9999'
$ws.Range("E45").Value = '  -0.07%  '

$ws.Range("B46").Value = 'TheSandbox'
$ws.Range("C46").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D46").Value = '0.4220'
$ws.Range("E46").Value = '  -0.53%  '

$ws.Range("D47").Value = '989.64'
$ws.Range("E47").Value = '  +2.15%  '

$ws.Range("D48").Value = '7.290'
$ws.Range("E48").Value = '  -1.84%  '

$ws.Range("D49").Value = '9.288'
$ws.Range("E49").Value = '  +5.95%  '

$ws.Range("D50").Value = '0.1165'
$ws.Range("E50").Value = '  -2.70%  '

$ws.Range("D51").Value = '34.19'
$ws.Range("E51").Value = '  -1.75%  '

# Restore original column style/format now that the text values are committed.
$ws.Range("D2:D51").Style = $origStyle
